$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.1029686666666667
$ws.Range("H2").Value = 0.308906
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6187893333333333
$ws.Range("N2").Value = 1.856368
$ws.Range("O2").Value = 0.1765034355725207
$ws.Range("P2").Value = 0.1765034355725208
$ws.Range("Q2").Value = 0.06371591260088888
$ws.Range("R2").Value = 0.5734432134080001
$ws.Range("S2").Value = 0.1765034355725207
$ws.Range("T2").Value = 0.1765034355725208

# Row 3 (Target cluster: FAPs)
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.1029686666666667
$ws.Range("H3").Value = 0.308906
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.233983666666667
$ws.Range("N3").Value = 6.701951000000001
$ws.Range("O3").Value = 0.6372213788099619
$ws.Range("P3").Value = 0.6372213788099619
$ws.Range("Q3").Value = 0.2300303195117778
$ws.Range("R3").Value = 2.070272875606
$ws.Range("S3").Value = 0.6372213788099619
$ws.Range("T3").Value = 0.6372213788099619

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.1029686666666667
$ws.Range("H4").Value = 0.308906
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.207158
$ws.Range("N4").Value = 0.6214740000000001
$ws.Range("O4").Value = 0.05908973658186135
$ws.Range("P4").Value = 0.05908973658186135
$ws.Range("Q4").Value = 0.02133078304933334
$ws.Range("R4").Value = 0.191977047444
$ws.Range("S4").Value = 0.05908973658186135
$ws.Range("T4").Value = 0.05908973658186135

# Row 5 (Target cluster: MuSCs)
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.1029686666666667
$ws.Range("H5").Value = 0.308906
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1889926666666667
$ws.Range("N5").Value = 0.566978
$ws.Range("O5").Value = 0.05390825789608347
$ws.Range("P5").Value = 0.05390825789608347
$ws.Range("Q5").Value = 0.01946032289644445
$ws.Range("R5").Value = 0.175142906068
$ws.Range("S5").Value = 0.05390825789608347
$ws.Range("T5").Value = 0.05390825789608347

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.1029686666666667
$ws.Range("H6").Value = 0.308906
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2568966666666667
$ws.Range("N6").Value = 0.77069
$ws.Range("O6").Value = 0.07327719113957255
$ws.Range("P6").Value = 0.07327719113957255
$ws.Range("Q6").Value = 0.02645230723777778
$ws.Range("R6").Value = 0.23807076514
$ws.Range("S6").Value = 0.07327719113957255
$ws.Range("T6").Value = 0.07327719113957255
